# Apply cryptocurrency price/volume updates as described by the commit diff.
# Two coin pairs (rows 27/28, 45/46, 49/50) were reordered (swapped) in the source ranking,
# and numeric D-column price cells must stay as exact text (e.g. trailing zeros kept),
# so we force the text number format ('@') on those cells before assigning their values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.129.11"
$ws.Range("E2").Value = "  +3.11%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.411.19"
$ws.Range("E3").Value = "  +4.12%  "

$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "406.89"
$ws.Range("E5").Value = "  -0.11%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.29"
$ws.Range("E6").Value = "  +17.84%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.610"
$ws.Range("E7").Value = "  +8.02%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.677"
$ws.Range("E9").Value = "  +10.08%  "

$ws.Range("E10").Value = "  +12.68%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "42.18"
$ws.Range("E11").Value = "  +9.44%  "

$ws.Range("E12").Value = "  -0.44%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.968.21"
$ws.Range("E13").Value = "  +4.95%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.56"
$ws.Range("E14").Value = "  +5.40%  "

$ws.Range("E15").Value = "  +4.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.417.34"
$ws.Range("E16").Value = "  +4.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.067.88"
$ws.Range("E17").Value = "  +3.49%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.55"
$ws.Range("E18").Value = "  +9.48%  "

$ws.Range("E19").Value = "  +5.38%  "

$ws.Range("E20").Value = "  +18.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "3.27"
$ws.Range("E21").Value = "  +0.07%  "

$ws.Range("E22").Value = "  +13.44%  "

$ws.Range("E23").Value = "  +6.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "308.64"
$ws.Range("E24").Value = "  +4.48%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.16"
$ws.Range("E25").Value = "  +3.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.60"
$ws.Range("E26").Value = "  +15.65%  "

$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.70"
$ws.Range("E27").Value = "  +10.40%  "

$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "29.80"
$ws.Range("E28").Value = "  +3.32%  "

$ws.Range("E29").Value = "  +1.83%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.47"
$ws.Range("E30").Value = "  +1.09%  "

$ws.Range("E31").Value = "  +2.81%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.81"
$ws.Range("E32").Value = "  +6.26%  "

$ws.Range("E33").Value = "  +7.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "42.68"
$ws.Range("E34").Value = "  +10.22%  "

$ws.Range("E35").Value = "  -0.07%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0489"
$ws.Range("E36").Value = "  +2.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "52.41"
$ws.Range("E37").Value = "  +0.65%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("E38").Value = "  +0.32%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.44"
$ws.Range("E39").Value = "  +4.49%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.98"
$ws.Range("E40").Value = "  -3.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.04"
$ws.Range("E41").Value = "  +9.94%  "

$ws.Range("E42").Value = "  +5.34%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "137.62"
$ws.Range("E43").Value = "  +2.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.99"
$ws.Range("E44").Value = "  +6.40%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.287"
$ws.Range("E45").Value = "  -1.70%  "

$ws.Range("B46").Value = "Celestia"
$ws.Range("C46").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.08"
$ws.Range("E46").Value = "  +5.67%  "

$ws.Range("E47").Value = "  +2.10%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.88"
$ws.Range("E48").Value = "  +5.00%  "

$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.749.44"
$ws.Range("E49").Value = "  +4.24%  "

$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.149.41"
$ws.Range("E50").Value = "  +1.89%  "

$ws.Range("E51").Value = "  -0.78%  "
